$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at the top of the data block (row 2). ---
# This shifts the old rows 2,3,4 down to 3,4,5, carrying values, styles
# and formulas with them (the C-column formula is re-pointed at its new
# row automatically).
$ws.Rows(2).Insert()

# --- 2. Remove all existing hyperlinks. ---
# Row-insert does not shift the <hyperlinks> refs, so the old entries
# still point at the pre-shift cells. Drop them all and re-add fresh
# ones at the correct, post-shift locations below.
$ws.Hyperlinks.Delete()

# --- 3. Populate the new row 2 (the Hampton Court fire entry). ---
# Borrow number/cell formatting from the equivalent cells one row down
# (which still carry the original style indexes after the insert), then
# write the real values on top.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("A2").Value = 31502
$ws.Range("E2").Value = "The Hampton Court fire. March 1986. – A retired London Fireman (home.blog)"

# --- 4. Re-create the six hyperlinks at their final (post-shift) cells.
# Passing the raw URL as the display text makes the engine record it as
# the <hyperlink display="..."> fallback (matching the original file's
# pattern of displaying the bare URL there); it also happens to overwrite
# the cell's visible text, so immediately restore each cell's real text
# afterwards. ---
$txt = $ws.Range("E3").Value()
$ws.Hyperlinks.Add($ws.Range("E3"), "http://news.bbc.co.uk/1/hi/england/london/6675381.stm", "", "", "http://news.bbc.co.uk/1/hi/england/london/6675381.stm") | Out-Null
$ws.Range("E3").Value = $txt

$txt = $ws.Range("F3").Value()
$ws.Hyperlinks.Add($ws.Range("F3"), "http://news.bbc.co.uk/1/hi/england/london/7643420.stm", "", "", "http://news.bbc.co.uk/1/hi/england/london/7643420.stm") | Out-Null
$ws.Range("F3").Value = $txt

$txt = $ws.Range("E5").Value()
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.thesun.co.uk/news/14840746/firefighters-tackle-huge-blaze-at-historic-boat-house/", "", "", "https://www.thesun.co.uk/news/14840746/firefighters-tackle-huge-blaze-at-historic-boat-house/") | Out-Null
$ws.Range("E5").Value = $txt

$txt = $ws.Range("D5").Value()
$ws.Hyperlinks.Add($ws.Range("D5"), "https://twitter.com/LondonFire/status/1389282224588132353?s=20") | Out-Null
$ws.Range("D5").Value = $txt

$txt = $ws.Range("E4").Value()
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.thesun.co.uk/news/9892888/worcester-park-fire-cause-sherbrooke-hamptons-london/", "", "", "https://www.thesun.co.uk/news/9892888/worcester-park-fire-cause-sherbrooke-hamptons-london/") | Out-Null
$ws.Range("E4").Value = $txt

$txt = $ws.Range("E2").Value()
$ws.Hyperlinks.Add($ws.Range("E2"), "https://beyondtheflamesandmore.home.blog/2021/03/31/the-hampton-court-fire-march-1986/", "", "", "https://beyondtheflamesandmore.home.blog/2021/03/31/the-hampton-court-fire-march-1986/") | Out-Null
$ws.Range("E2").Value = $txt

# --- 5. `Hyperlinks.Add` stamps a brand-new "hyperlink look" style on the
# cell it touches instead of reusing the workbook's existing Hyperlink
# style index; restore the original per-cell formatting (copied from the
# untouched, still-correctly-styled F4 hyperlink cell) without disturbing
# the links or the text restored above. ---
$ws.Range("F4").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)

# --- 6. Restore the selection shown in the saved view. ---
$ws.Range("E10").Select() | Out-Null
